# "Commit connected all the program"
#
# - Sheet3 gets a new third row: A3 = "abc"
# - Sheet3 becomes the active sheet/tab (selection moves to A3), so
#   Sheet1 is no longer the tab-selected sheet.

$wb = $excel.ActiveWorkbook

$ws3 = $wb.Worksheets.Item("Sheet3")

# Add the new value to Sheet3
$ws3.Range("A3").Value = "abc"

# Make Sheet3 the active sheet and select the new cell,
# which moves tabSelected from Sheet1 to Sheet3.
$ws3.Activate()
$ws3.Range("A3").Select()
